$wb = $excel.ActiveWorkbook

# --- 1. Create new sheet "2.1" at the end of the workbook, based on "1.4" ---
$src = $wb.Worksheets.Item("1.4")
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$src.Copy($null, $last)
$new = $wb.Worksheets.Item($wb.Worksheets.Count)
$new.Name = "2.1"

# Remove the copied chart/drawing - sheet 2.1 has no chart
foreach ($co in $new.ChartObjects()) {
    $co.Delete()
}

# --- 2. Edit the new sheet content ---
# Header row (B2) text is unchanged: "Burnt Down Chart Sprint 1 - Hari Ketiga"

# Row 6 headers: Jumat, Sabtu, Minggu instead of Senin..Kamis; drop column H
$new.Range("E6").Value = "Jumat"
$new.Range("F6").Value = "Sabtu"
$new.Range("G6").Value = "Minggu"

# Backlog rows 7-10
$new.Range("B7").Value = "Question Tag"
$new.Range("B8").Value = "Conditional Sentence"
$new.Range("B9").Value = "Many & Much"
$new.Range("B10").Value = "Preposition"

$new.Range("E7").Value = 100
$new.Range("F7").Value = "-"
$new.Range("G7").Value = "-"

$new.Range("E8").Value = 100
$new.Range("F8").Value = "-"
$new.Range("G8").Value = "-"

$new.Range("E9").Value = "-"
$new.Range("F9").Value = "-"
$new.Range("G9").Value = "-"

$new.Range("E10").Value = "-"
$new.Range("F10").Value = "-"
$new.Range("G10").Value = "-"

# Remove old row 11 backlog entry (Tenses IV) and shift totals row up to row 11
$new.Range("B11:H11").Delete()

$new.Range("D11").Value = "Total = "
$new.Range("E11").Formula = "=SUM(E7:E10)"
$new.Range("F11").Formula = "=SUM(F7:F10)"
$new.Range("G11").Formula = "=SUM(G7:G10)"

# Drop column H entirely (only B:G used now)
$new.Columns.Item("H").Delete()

$new.Range("E11").Select()

Write-Host "Done"
